# Added Working Excel Deserializer
# This script re-applies the flight-record rows that the deserializer wrote
# to Sheet1: it fixes up row 2's number formatting (scheduled time / seats /
# year-of-registration) and appends four more flight rows (3-6) that follow
# the same shape, then leaves the selection where the deserializer left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlLeft = -4131

function Set-FlightRow {
    param(
        $Row,
        $ScheduledTime,
        $DestCode,
        $DestName,
        $FlightType,
        $Manufacturer,
        $Model,
        $Seats,
        $RegNumber,
        $YearOfReg,
        $PlaneState,
        $Airline,
        $Terminal
    )

    # Text columns: plain left-aligned general format.
    $ws.Cells.Item($Row, 2).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 2).Value = $DestCode

    $ws.Cells.Item($Row, 3).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 3).Value = $DestName

    $ws.Cells.Item($Row, 4).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 4).Value = $FlightType

    $ws.Cells.Item($Row, 5).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 5).Value = $Manufacturer

    $ws.Cells.Item($Row, 6).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 6).Value = $Model

    $ws.Cells.Item($Row, 8).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 8).Value = $RegNumber

    $ws.Cells.Item($Row, 10).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 10).Value = $PlaneState

    $ws.Cells.Item($Row, 11).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 11).Value = $Airline

    $ws.Cells.Item($Row, 12).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 12).Value = $Terminal

    # Column A: scheduled time - stored as literal text but carries a date
    # display format (mm-dd-yy == built-in numFmtId 14).
    $ws.Cells.Item($Row, 1).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 1).NumberFormat = "mm-dd-yy"
    $ws.Cells.Item($Row, 1).Value = $ScheduledTime

    # Columns G/I: numeric columns formatted as plain integers (numFmtId 1).
    $ws.Cells.Item($Row, 7).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 7).NumberFormat = "0"
    $ws.Cells.Item($Row, 7).Value = $Seats

    $ws.Cells.Item($Row, 9).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($Row, 9).NumberFormat = "0"
    $ws.Cells.Item($Row, 9).Value = $YearOfReg
}

# Row 2 already holds its flight record - just re-apply it so the
# number-format styles (date / integer) get created and attached.
Set-FlightRow 2 "2017-06-07T13:34:08.0039447-05:00" `
    "LBBS" "Burgas Airport" "Arrivals" `
    "Embraer" "E175" 80 "BG66666SL" `
    2005 "Good condition" "Hemus Air" "A5"

Set-FlightRow 3 "2017-06-07T13:34:08.0039447-05:01" `
    "LBBS" "Burgas Airport" "Arrivals" `
    "Embraer" "E176" 81 "BG66666SL" `
    2006 "Good condition" "Hemus Air" "A6"

Set-FlightRow 4 "2017-06-07T13:34:08.0039447-05:02" `
    "LBBS" "Burgas Airport" "Arrivals" `
    "Embraer" "E177" 82 "BG66666SL" `
    2007 "Good condition" "Hemus Air" "A7"

Set-FlightRow 5 "2017-06-07T13:34:08.0039447-05:03" `
    "LBBS" "Burgas Airport" "Arrivals" `
    "Embraer" "E178" 83 "BG66666SL" `
    2008 "Good condition" "Hemus Air" "A8"

Set-FlightRow 6 "2017-06-07T13:34:08.0039447-05:04" `
    "LBBS" "Burgas Airport" "Arrivals" `
    "Embraer" "E179" 84 "BG66666SL" `
    2009 "Good condition" "Hemus Air" "A9"

# Leave the selection where the deserializer run left it.
$ws.Range("I9").Select()
